$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.924.42'
$ws.Range('E2').Value = '  +1.29%  '
$ws.Range('D3').Value = '2.707.18'
$ws.Range('E3').Value = '  +2.69%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = "'609.10"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.09%  '
$ws.Range('D6').Value = "'158.16"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.23%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = "'0.588"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.56%  '
$ws.Range('E9').Value = '  +5.50%  '
$ws.Range('E10').Value = '  +3.86%  '
$ws.Range('E11').Value = '  +0.40%  '
$ws.Range('E12').Value = '  +1.04%  '
$ws.Range('D13').Value = "'30.49"
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').Value = "'0.0000203"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +8.93%  '
$ws.Range('D15').Value = '3.196.46'
$ws.Range('E15').Value = '  +2.73%  '
$ws.Range('D16').Value = '65.808.61'
$ws.Range('E16').Value = '  +1.25%  '
$ws.Range('D17').Value = '2.706.46'
$ws.Range('E17').Value = '  +1.37%  '
$ws.Range('D18').Value = "'12.73"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.24%  '
$ws.Range('E19').Value = '  +1.98%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').Value = "'360.00"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.85%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value = "'7.66"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.26%  '
$ws.Range('E22').Value = '  -0.07%  '
$ws.Range('D23').Value = "'70.93"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.85%  '
$ws.Range('D24').Value = "'9.89"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E25').Value = '  +12.61%  '
$ws.Range('E26').Value = '  -1.69%  '
$ws.Range('E27').Value = '  +3.43%  '
$ws.Range('E28').Value = '  +3.88%  '
$ws.Range('D29').Value = "'8.43"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.25%  '
$ws.Range('E30').Value = '  +5.09%  '
$ws.Range('D31').Value = "'544.29"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.43%  '
$ws.Range('E32').Value = '  +0.03%  '
$ws.Range('E33').Value = '  +2.44%  '
$ws.Range('D34').Value = "'6.77"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +6.60%  '
$ws.Range('D35').Value = "'5.44"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.73%  '
$ws.Range('E36').Value = '  +2.25%  '
$ws.Range('D37').Value = "'20.94"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.02%  '
$ws.Range('D38').Value = "'163.35"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.07%  '
$ws.Range('E39').Value = '  -0.26%  '
$ws.Range('D40').Value = "'1.00"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.02%  '
$ws.Range('D41').Value = "'173.08"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.39%  '
$ws.Range('E42').Value = '  +0.03%  '
$ws.Range('E43').Value = '  +0.68%  '
$ws.Range('E44').Value = '  +2.92%  '
$ws.Range('E45').Value = '  +0.33%  '
$ws.Range('D46').Value = "'23.57"
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').Value = "'2.29"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.89%  '
$ws.Range('E48').Value = '  +4.35%  '
$ws.Range('E49').Value = '  +1.46%  '
$ws.Range('D50').Value = "'21.08"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +8.14%  '
$ws.Range('D51').Value = "'0.0993"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.03%  '
